$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "G2"  = 1.73
    "H2"  = 3.5
    "I2"  = 5.25
    "J2"  = 2.4
    "K2"  = 2
    "L2"  = 6
    "M2"  = 1.08
    "N2"  = 8
    "O2"  = 1.44
    "P2"  = 2.63
    "Q2"  = 2.4
    "R2"  = 1.53
    "S2"  = 1.53
    "T2"  = 2.38
    "X2"  = 7
    "Y2"  = 9
    "Z2"  = 13
    "AA2" = 17
    "AD2" = 7
    "AG2" = 10
    "AH2" = 23
    "AN2" = 3.5
    "AO2" = 9.5
    "AQ2" = 34
    "AT2" = 2.38
    "AW2" = 6.5
    "AX2" = 34
    "AZ2" = 126
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
